# Case_4_11/lines_states.xlsx update: add line7/line8 rows (pushing the
# extr1..extr8 rows down by two) and tweak several from_bus/to_bus/in_service
# values, per the "contingencies with rene fine" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of an existing data row (row 8, A:E) into the two
# brand-new rows (16 and 17) so they pick up the same bold/bordered/centered
# style used by every other data row, instead of Excel inventing a fresh
# (duplicate) style for them.
$ws.Range("A8:E8").Copy() | Out-Null
$ws.Range("A16:E16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A8:E8").Copy() | Out-Null
$ws.Range("A17:E17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Final state for rows 8-17 (name/from_bus/to_bus/in_service), with the new
# "line7"/"line8" rows inserted right after "line6" and the old extr1-8 rows
# shifted down two rows to extr1-8 @ rows 10-17.
$rows = @(
    @{ R=8;  A=6;  B="line7"; C=14; D=11; E=$true  },
    @{ R=9;  A=7;  B="line8"; C=16; D=9;  E=$true  },
    @{ R=10; A=8;  B="extr1"; C=5;  D=12; E=$true  },
    @{ R=11; A=9;  B="extr2"; C=5;  D=9;  E=$true  },
    @{ R=12; A=10; B="extr3"; C=10; D=11; E=$true  },
    @{ R=13; A=11; B="extr4"; C=7;  D=8;  E=$true  },
    @{ R=14; A=12; B="extr5"; C=9;  D=11; E=$false },
    @{ R=15; A=13; B="extr6"; C=7;  D=11; E=$true  },
    @{ R=16; A=14; B="extr7"; C=5;  D=7;  E=$true  },
    @{ R=17; A=15; B="extr8"; C=8;  D=5;  E=$false }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value2 = $row.A
    $ws.Cells.Item($r, 2).Value2 = $row.B
    $ws.Cells.Item($r, 3).Value2 = $row.C
    $ws.Cells.Item($r, 4).Value2 = $row.D
    $ws.Cells.Item($r, 5).Value2 = $row.E
}
